# Apply crypto price/volume updates scraped on Sun Oct 22 04:49:07 UTC 2023
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.968.10"
$ws.Range("E2").Value = "  +1.19%  "
$ws.Range("E3").Value = "  +2.06%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.07"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.35%  "
$ws.Range("E6").Value = "  +1.26%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "29.67"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +10.49%  "
$ws.Range("E9").Value = "  +3.86%  "
$ws.Range("E10").Value = "  +2.29%  "
$ws.Range("E11").Value = "  +0.73%  "
$ws.Range("D12").Value = "1.867.47"
$ws.Range("E12").Value = "  +2.13%  "
$ws.Range("D13").Value = "1.635.86"
$ws.Range("E13").Value = "  +2.11%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.574"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +6.77%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "9.54"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +25.05%  "
$ws.Range("E16").Value = "  +4.25%  "
$ws.Range("D17").Value = "29.992.51"
$ws.Range("E17").Value = "  +1.30%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "64.91"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "249.44"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +3.28%  "
$ws.Range("E20").Value = "  +2.00%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.999"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.14%  "
$ws.Range("E22").Value = "  +5.02%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.68"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +4.97%  "
$ws.Range("E24").Value = "  +1.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "160.06"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +3.30%  "
$ws.Range("E26").Value = "  +2.53%  "
$ws.Range("E27").Value = "  +2.78%  "
$ws.Range("E28").Value = "  +3.77%  "
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("E30").Value = "  +2.79%  "
$ws.Range("E31").Value = "  +6.25%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.40"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +5.38%  "
$ws.Range("E33").Value = "  +2.30%  "
$ws.Range("D34").Value = "1.434.81"
$ws.Range("E34").Value = "  +0.70%  "
$ws.Range("E35").Value = "  +7.93%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.04"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.87"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.37%  "
$ws.Range("E38").Value = "  +3.00%  "
$ws.Range("E39").Value = "  -0.25%  "
$ws.Range("E40").Value = "  +2.67%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "73.66"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +11.17%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.833"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +3.14%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.98"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.86%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "54.91"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +1.72%  "
$ws.Range("E45").Value = "  +0.21%  "
$ws.Range("E46").Value = "  +6.12%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.11%  "
$ws.Range("E48").Value = "  +3.46%  "
$ws.Range("D49").Value = "1.773.84"
$ws.Range("E49").Value = "  +2.01%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "90.43"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +5.15%  "
$ws.Range("E51").Value = "  +3.24%  "

Write-Output "Updated 75 cells"
